# The three pictures that live in the document's headers/footers (the
# Pearson logo used on the "default" and "first page" footers, and the
# BTec logo used on the "first page" header) each carry a display Name
# in two places inside their drawing markup: <wp:docPr .../> and the
# nested <pic:cNvPr .../>. Real Word keeps those two attributes in sync
# but exposes no InlineShape.Name property on the object model to edit
# them directly, so we round-trip the package's WordOpenXML and patch
# the three <wp:docPr>/<pic:cNvPr> name pairs in place - each scoped by
# its picture's unique wp:docPr id so the three pictures can't
# cross-match each other.

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

$opts = [System.Text.RegularExpressions.RegexOptions]::Singleline

# Header "first page" (header1.xml): BTec_Logo-Orange, wp:docPr id="1"
$pattern1 = '(<wp:docPr\b[^>]*\bid="1"[^>]*\bname=")image1\.jpg("[^>]*/>.*?<pic:cNvPr\b[^>]*\bname=")image1\.jpg("[^>]*/>)'
$repl1 = '${1}image2.jpg${2}image2.jpg${3}'
$xml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern1, $repl1, $opts)

# Footer "default" (footer2.xml): PearsonLogo, wp:docPr id="2"
$pattern2 = '(<wp:docPr\b[^>]*\bid="2"[^>]*\bname=")image2\.png("[^>]*/>.*?<pic:cNvPr\b[^>]*\bname=")image2\.png("[^>]*/>)'
$repl2 = '${1}image1.png${2}image1.png${3}'
$xml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern2, $repl2, $opts)

# Footer "first page" (footer1.xml): PearsonLogo, wp:docPr id="3"
$pattern3 = '(<wp:docPr\b[^>]*\bid="3"[^>]*\bname=")image2\.png("[^>]*/>.*?<pic:cNvPr\b[^>]*\bname=")image2\.png("[^>]*/>)'
$repl3 = '${1}image1.png${2}image1.png${3}'
$xml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern3, $repl3, $opts)

$d.Content.WordOpenXML = $xml

Write-Output "done"
